$wb = $excel.ActiveWorkbook

# Rename the "High_School_Units" sheet to "High School Units"
$ws = $wb.Worksheets.Item("High_School_Units")
$ws.Name = "High School Units"

# Fix the value in A7 (was 3, should be 4)
$ws.Range("A7").Value = 4

# Update the active selection on that sheet to A8
$ws.Activate()
$ws.Range("A8").Select()
